# ===========================================================================
# LMS2 storage type - update data values (2017 -> 2022 column reshuffle)
# Sheet 1 "weighted values by province": rows 2-21, columns F:Q updated.
# Sheet 2 "wetight percentage by province": rows 2 & 12 year headers updated,
#   plus the active-cell selection moved to C2:R21.
# ===========================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "weighted values by province"
$ws2 = $wb.Worksheets.Item(2)   # "wetight percentage by province"

# --- Sheet 1: rewrite the F:Q block for every affected row ----------------
$row2 = New-Object "object[,]" 1,12
$row2[0,0] = 2022
$row2[0,1] = 2017
$row2[0,2] = 2022
$row2[0,3] = 2017
$row2[0,4] = 2022
$row2[0,5] = 2017
$row2[0,6] = 2022
$row2[0,7] = 2017
$row2[0,8] = 2022
$row2[0,9] = 2017
$row2[0,10] = 2022
$row2[0,11] = 2017
$ws1.Range("F2:Q2").Value = $row2

$row3 = New-Object "object[,]" 1,12
$row3[0,0] = 8.8000000000000007
$row3[0,1] = 162.1
$row3[0,2] = 172
$row3[0,3] = 8.5
$row3[0,4] = 13.4
$row3[0,5] = 16.399999999999999
$row3[0,6] = 4.0999999999999996
$row3[0,7] = 4.9000000000000004
$row3[0,8] = 0
$row3[0,9] = 2.2000000000000002
$row3[0,10] = 7.6
$row3[0,11] = 5.0999999999999996
$ws1.Range("F3:Q3").Value = $row3

$row4 = New-Object "object[,]" 1,12
$row4[0,0] = 13.2
$row4[0,1] = 61.3
$row4[0,2] = 47.1
$row4[0,3] = 44.1
$row4[0,4] = 45.5
$row4[0,5] = 17.7
$row4[0,6] = 13.6
$row4[0,7] = 59.2
$row4[0,8] = 47.9
$row4[0,9] = 27.5
$row4[0,10] = 22.7
$row4[0,11] = 2.2000000000000002
$ws1.Range("F4:Q4").Value = $row4

$row5 = New-Object "object[,]" 1,12
$row5[0,0] = 0
$row5[0,1] = 45.7
$row5[0,2] = 44.4
$row5[0,3] = 0
$row5[0,4] = 6
$row5[0,5] = 2.2000000000000002
$row5[0,6] = 5.6
$row5[0,7] = 0
$row5[0,8] = 4.2
$row5[0,9] = 0
$row5[0,10] = 4.3
$row5[0,11] = 5.0999999999999996
$ws1.Range("F5:Q5").Value = $row5

$row6 = New-Object "object[,]" 1,12
$row6[0,0] = 0
$row6[0,1] = 57.1
$row6[0,2] = 0
$row6[0,3] = 12.4
$row6[0,4] = 0
$row6[0,5] = 6.9
$row6[0,6] = 0
$row6[0,7] = 13.7
$row6[0,8] = 0
$row6[0,9] = 7.5
$row6[0,10] = 0
$row6[0,11] = 3.6
$ws1.Range("F6:Q6").Value = $row6

$row7 = New-Object "object[,]" 1,12
$row7[0,0] = 0
$row7[0,1] = 40.4
$row7[0,2] = 0
$row7[0,3] = 24
$row7[0,4] = 0
$row7[0,5] = 6.1
$row7[0,6] = 0
$row7[0,7] = 21.2
$row7[0,8] = 0
$row7[0,9] = 25.5
$row7[0,10] = 0
$row7[0,11] = 10.199999999999999
$ws1.Range("F7:Q7").Value = $row7

$row8 = New-Object "object[,]" 1,12
$row8[0,0] = 23
$row8[0,1] = 424.9
$row8[0,2] = 488.8
$row8[0,3] = 342.2
$row8[0,4] = 311.7
$row8[0,5] = 234
$row8[0,6] = 139.5
$row8[0,7] = 615.20000000000005
$row8[0,8] = 688
$row8[0,9] = 193
$row8[0,10] = 174.7
$row8[0,11] = 62.1
$ws1.Range("F8:Q8").Value = $row8

$row9 = New-Object "object[,]" 1,12
$row9[0,0] = 59.3
$row9[0,1] = 573.5
$row9[0,2] = 400.1
$row9[0,3] = 38.200000000000003
$row9[0,4] = 133.19999999999999
$row9[0,5] = 438.8
$row9[0,6] = 562.29999999999995
$row9[0,7] = 1687.4
$row9[0,8] = 1443.7
$row9[0,9] = 16.899999999999999
$row9[0,10] = 63.5
$row9[0,11] = 64.599999999999994
$ws1.Range("F9:Q9").Value = $row9

$row10 = New-Object "object[,]" 1,12
$row10[0,0] = 0
$row10[0,1] = 0
$row10[0,2] = 41
$row10[0,3] = 0
$row10[0,4] = 0
$row10[0,5] = 0
$row10[0,6] = 0
$row10[0,7] = 0
$row10[0,8] = 1.9
$row10[0,9] = 0
$row10[0,10] = 0
$row10[0,11] = 0
$ws1.Range("F10:Q10").Value = $row10

$row11 = New-Object "object[,]" 1,12
$row11[0,0] = 104.3
$row11[0,1] = 1365
$row11[0,2] = 1193.4000000000001
$row11[0,3] = 469.4
$row11[0,4] = 509.8
$row11[0,5] = 722.1
$row11[0,6] = 725.1
$row11[0,7] = 2401.6
$row11[0,8] = 2185.6999999999998
$row11[0,9] = 272.60000000000002
$row11[0,10] = 272.8
$row11[0,11] = 152.9
$ws1.Range("F11:Q11").Value = $row11

$row12 = New-Object "object[,]" 1,12
$row12[0,0] = 2022
$row12[0,1] = 2017
$row12[0,2] = 2022
$row12[0,3] = 2017
$row12[0,4] = 2022
$row12[0,5] = 2017
$row12[0,6] = 2022
$row12[0,7] = 2017
$row12[0,8] = 2022
$row12[0,9] = 2017
$row12[0,10] = 2022
$row12[0,11] = 2017
$ws1.Range("F12:Q12").Value = $row12

$row13 = New-Object "object[,]" 1,12
$row13[0,0] = 6.7
$row13[0,1] = 63.3
$row13[0,2] = 43.5
$row13[0,3] = 0
$row13[0,4] = 7.5
$row13[0,5] = 6.5
$row13[0,6] = 2.9
$row13[0,7] = 0
$row13[0,8] = 0
$row13[0,9] = 0
$row13[0,10] = 3.3
$row13[0,11] = 6.2
$ws1.Range("F13:Q13").Value = $row13

$row15 = New-Object "object[,]" 1,12
$row15[0,0] = 27.2
$row15[0,1] = 74.5
$row15[0,2] = 61
$row15[0,3] = 0
$row15[0,4] = 0
$row15[0,5] = 23.2
$row15[0,6] = 4.4000000000000004
$row15[0,7] = 4.4000000000000004
$row15[0,8] = 4.7
$row15[0,9] = 0
$row15[0,10] = 0
$row15[0,11] = 0
$ws1.Range("F15:Q15").Value = $row15

$row18 = New-Object "object[,]" 1,12
$row18[0,0] = 50.4
$row18[0,1] = 51.6
$row18[0,2] = 34.299999999999997
$row18[0,3] = 88.4
$row18[0,4] = 46.5
$row18[0,5] = 53.5
$row18[0,6] = 36.9
$row18[0,7] = 113.6
$row18[0,8] = 85.1
$row18[0,9] = 338
$row18[0,10] = 380.1
$row18[0,11] = 4.0999999999999996
$ws1.Range("F18:Q18").Value = $row18

$row19 = New-Object "object[,]" 1,12
$row19[0,0] = 30.5
$row19[0,1] = 54.3
$row19[0,2] = 50.2
$row19[0,3] = 27.3
$row19[0,4] = 15.3
$row19[0,5] = 202.7
$row19[0,6] = 193.5
$row19[0,7] = 525.1
$row19[0,8] = 405.2
$row19[0,9] = 15.8
$row19[0,10] = 28.3
$row19[0,11] = 3.8
$ws1.Range("F19:Q19").Value = $row19

$row20 = New-Object "object[,]" 1,12
$row20[0,0] = 0
$row20[0,1] = 0
$row20[0,2] = 6.5
$row20[0,3] = 0
$row20[0,4] = 0
$row20[0,5] = 0
$row20[0,6] = 0
$row20[0,7] = 0
$row20[0,8] = 0
$row20[0,9] = 0
$row20[0,10] = 0
$row20[0,11] = 0
$ws1.Range("F20:Q20").Value = $row20

$row21 = New-Object "object[,]" 1,12
$row21[0,0] = 114.8
$row21[0,1] = 243.7
$row21[0,2] = 195.5
$row21[0,3] = 115.7
$row21[0,4] = 69.3
$row21[0,5] = 285.89999999999998
$row21[0,6] = 237.7
$row21[0,7] = 643.1
$row21[0,8] = 495
$row21[0,9] = 353.8
$row21[0,10] = 411.7
$row21[0,11] = 14.1
$ws1.Range("F21:Q21").Value = $row21

# --- Sheet 2: year headers on rows 2 and 12 --------------------------------
$hdr = New-Object "object[,]" 1,6
$hdr[0,0] = 2022   # F
$hdr[0,1] = 2022   # H
$hdr[0,2] = 2022   # J
$hdr[0,3] = 2017   # M
$hdr[0,4] = 2017   # O
$hdr[0,5] = 2017   # Q

foreach ($r in 2, 12) {
    $ws2.Range("F$r").Value = $hdr[0,0]
    $ws2.Range("H$r").Value = $hdr[0,1]
    $ws2.Range("J$r").Value = $hdr[0,2]
    $ws2.Range("M$r").Value = $hdr[0,3]
    $ws2.Range("O$r").Value = $hdr[0,4]
    $ws2.Range("Q$r").Value = $hdr[0,5]
}

# --- Sheet 2: move the saved selection/active cell to C2:R21 --------------
$ws2.Activate()
$ws2.Range("C2:R21").Select()

"edit complete"
